$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "filas" and "filasE": add new row 230 for the new "GTMva000" (Valor
#    agregado) row, belonging to area 5.
# ---------------------------------------------------------------------------
$filas = $wb.Worksheets.Item("filas")
$filas.Range("A230").Value = "GTMva000"
$filas.Range("B230").Value = 5

$filasE = $wb.Worksheets.Item("filasE")
$filasE.Range("A230").Value = "GTMva000"
$filasE.Range("B230").Value = 5

# ---------------------------------------------------------------------------
# 2) "areas_filas": add new row 7 describing area 5 = "Valor agregado"
# ---------------------------------------------------------------------------
$areasFilas = $wb.Worksheets.Item("areas_filas")
$areasFilas.Range("A7").Value = 5
$areasFilas.Range("B7").Value = "Valor agregado"

# ---------------------------------------------------------------------------
# 3) "cuadros": add new row 11 describing cuadro 10 = "Valor Agregado"
# ---------------------------------------------------------------------------
$cuadros = $wb.Worksheets.Item("cuadros")
$cuadros.Range("A11").Value = 10
$cuadros.Range("B11").Value = "Valor Agregado"

# ---------------------------------------------------------------------------
# 4) "npg": fix two existing text labels
# ---------------------------------------------------------------------------
$npg = $wb.Worksheets.Item("npg")
$npg.Range("B143").Value = "Energía eléctrica, gas, vapor y aire acondicionado"
$npg.Range("B221").Value = "Otras primarias (incluye solar)"

# ---------------------------------------------------------------------------
# 5) "ciiu": add compact classification columns D (id_ciiu1_compacta) and
#    E (ciiu1_compacta)
# ---------------------------------------------------------------------------
$ciiu = $wb.Worksheets.Item("ciiu")
$ciiu.Range("D1").Value = "id_ciiu1_compacta"
$ciiu.Range("E1").Value = "ciiu1_compacta"

# Write the cells that introduce brand-new shared strings first, in the same
# order the original author typed them, so the shared string table ends up
# in the same order as the target workbook.
$ciiu.Range("E22").Value = "No determinada"
$ciiu.Range("D10").Value = "I-S"
$ciiu.Range("E9").Value = "Transporte"
$ciiu.Range("D5").Value = "D-E"

$ciiu.Range("D2").Value = "A"
$ciiu.Range("E2").Value = "Agricultura"

$ciiu.Range("D3").Value = "B"
$ciiu.Range("E3").Value = "Minería"

$ciiu.Range("D4").Value = "C"
$ciiu.Range("E4").Value = "Manufacturas"

$ciiu.Range("E5").Value = "Servicios básicos"

$ciiu.Range("D6").Value = "D-E"
$ciiu.Range("E6").Value = "Servicios básicos"

$ciiu.Range("D7").Value = "F"
$ciiu.Range("E7").Value = "Construcción"

$ciiu.Range("D8").Value = "G"
$ciiu.Range("E8").Value = "Comercio"

$ciiu.Range("D9").Value = "H"

$ciiu.Range("E10").Value = "Otros servicios"

$ciiu.Range("D11").Value = "I-S"
$ciiu.Range("E11").Value = "Otros servicios"

$ciiu.Range("D12").Value = "I-S"
$ciiu.Range("E12").Value = "Otros servicios"

$ciiu.Range("D13").Value = "I-S"
$ciiu.Range("E13").Value = "Otros servicios"

$ciiu.Range("D14").Value = "I-S"
$ciiu.Range("E14").Value = "Otros servicios"

$ciiu.Range("D15").Value = "I-S"
$ciiu.Range("E15").Value = "Otros servicios"

$ciiu.Range("D16").Value = "I-S"
$ciiu.Range("E16").Value = "Otros servicios"

$ciiu.Range("D17").Value = "I-S"
$ciiu.Range("E17").Value = "Otros servicios"

$ciiu.Range("D18").Value = "I-S"
$ciiu.Range("E18").Value = "Otros servicios"

$ciiu.Range("D19").Value = "I-S"
$ciiu.Range("E19").Value = "Otros servicios"

$ciiu.Range("D20").Value = "I-S"
$ciiu.Range("E20").Value = "Otros servicios"

$ciiu.Range("D21").Value = "I-S"
$ciiu.Range("E21").Value = "Otros servicios"

$ciiu.Range("D22").Value = "Z"

# ---------------------------------------------------------------------------
# 6) "naeg": add compact classification columns F/G, looked up from "ciiu"
# ---------------------------------------------------------------------------
$naeg = $wb.Worksheets.Item("naeg")
$naeg.Range("F1").Value = "id_ciiu1_compacta"
$naeg.Range("G1").Value = "ciiu1_compacta"

for ($r = 2; $r -le 134; $r++) {
    $naeg.Range("F$r").Formula = "=VLOOKUP(C$r,ciiu!`$A`$2:`$E`$22,4,FALSE)"
    $naeg.Range("G$r").Formula = "=VLOOKUP(C$r,ciiu!`$A`$2:`$E`$22,5,FALSE)"
}

$wb.Save()
